$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue 'D2' '291.77'
Set-TextValue 'E2' '-7.98%'
Set-TextValue 'D3' '40.48'
Set-TextValue 'E3' '-1.48%'
Set-TextValue 'E4' '-2.51%'
Set-TextValue 'D5' '0.07290'
Set-TextValue 'E5' '-4.51%'
Set-TextValue 'B6' 'GateToken'
Set-TextValue 'C6' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue 'D6' '4.282'
Set-TextValue 'E6' '-0.86%'
Set-TextValue 'B7' 'FTXToken'
Set-TextValue 'C7' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue 'D7' '1.560'
Set-TextValue 'E7' '-7.08%'
Set-TextValue 'B8' 'MXToken'
Set-TextValue 'C8' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D8' '0.9230'
Set-TextValue 'E8' '-1.14%'
Set-TextValue 'B9' 'LiechtensteinCryptoassetsExchange'
Set-TextValue 'C9' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D9' '0.1158'
Set-TextValue 'E9' '-6.89%'
Set-TextValue 'B10' 'WazirX'
Set-TextValue 'C10' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D10' '0.1731'
Set-TextValue 'E10' '-5.18%'
Set-TextValue 'B11' 'MandalaExchangeToken'
Set-TextValue 'C11' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D11' '0.08579'
Set-TextValue 'E11' '-5.55%'
Set-TextValue 'B12' 'BitrueCoin'
Set-TextValue 'C12' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D12' '0.04178'
Set-TextValue 'E12' '0.88%'
Set-TextValue 'B13' 'BitMartToken'
Set-TextValue 'C13' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D13' '0.1053'
Set-TextValue 'E13' '-0.37%'
Set-TextValue 'B14' 'BitForexToken'
Set-TextValue 'C14' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D14' '0.001264'
Set-TextValue 'E14' '-2.09%'
Set-TextValue 'B15' 'TigerCash'
Set-TextValue 'C15' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue 'D15' '0.005961'
Set-TextValue 'E15' '0.47%'
Set-TextValue 'B16' 'LEO'
Set-TextValue 'C16' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D16' '3.397'
Set-TextValue 'E16' '1.37%'
Set-TextValue 'E18' '-2.50%'
Set-TextValue 'D19' '7.861'
Set-TextValue 'E19' '-6.51%'
Set-TextValue 'E20' '2.52%'
Set-TextValue 'D21' '0.2883'
Set-TextValue 'E21' '0.48%'
Set-TextValue 'D22' '0.03855'
Set-TextValue 'E22' '-4.66%'
Set-TextValue 'E23' '-0.96%'
Set-TextValue 'D24' '0.003793'
Set-TextValue 'E24' '-7.12%'
Set-TextValue 'D25' '0.0001280'
Set-TextValue 'E25' '0.39%'
Set-TextValue 'D26' '0.0003725'
Set-TextValue 'D38' '0.02305'
Set-TextValue 'E38' '-7.95%'
Set-TextValue 'D39' '0.04966'
Set-TextValue 'E39' '-5.35%'
Set-TextValue 'E40' '214.64%'
Set-TextValue 'D41' '0.007692'
Set-TextValue 'E41' '-1.19%'
Set-TextValue 'D42' '0.1273'
Set-TextValue 'E42' '-1.89%'
Set-TextValue 'D43' '0.007372'
Set-TextValue 'E43' '4.12%'
Set-TextValue 'D44' '0.007072'
Set-TextValue 'E44' '-14.17%'
Set-TextValue 'D45' '0.3150'
Set-TextValue 'E45' '-0.44%'
Set-TextValue 'D46' '0.00006422'
Set-TextValue 'E46' '-3.92%'
Set-TextValue 'D47' '0.00000000751'
Set-TextValue 'E47' '-0.29%'
Set-TextValue 'D48' '0.01504'
Set-TextValue 'E48' '-93.26%'
Set-TextValue 'E49' '-0.35%'
Set-TextValue 'D50' '0.00002102'
Set-TextValue 'E50' '-0.29%'
Set-TextValue 'E51' '-0.29%'
